# Fix header labels on the existing sheets and add a new "PO Forecast" sheet
# with forecasted purchase-order quantities (ds / PO_Forecast / yhat_lower / yhat_upper).

$wb = $excel.ActiveWorkbook

# --- Rename the generic "Requested quantity" headers to clearer field names ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after the last existing sheet ---
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

$headerRange = $wsForecast.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Forecast data rows
$data = @(
  @(45158.99999999999, 22, 3.511080157598952, 42.97893227888814),
  @(45214.99999999999, 21, 1.849369167887289, 40.16038137717783),
  @(45221.99999999999, 21, 2.011892318412262, 42.82916642283537),
  @(45228.99999999999, 21, 0.8945938009170535, 42.19090457472137),
  @(45235.99999999999, 21, 1.404622641176621, 42.14357232600431),
  @(45242.99999999999, 21, 0.4693034198461432, 40.18976909060613),
  @(45249.99999999999, 21, 1.538738414251331, 41.81139007650135),
  @(45256.99999999999, 21, 1.003493266263698, 41.36078325807129),
  @(45263.99999999999, 21, 0.4522213610179129, 40.65722866018895),
  @(45270.99999999999, 21, 0.1225132077828309, 42.24108404309036),
  @(45277.99999999999, 20, 1.655170271397111, 41.85051233427591),
  @(45298.99999999999, 20, -0.531481979124406, 40.68215166948372),
  @(45382.99999999999, 19, -1.068192615017379, 39.35477651614787),
  @(45389.99999999999, 19, -1.574493359785864, 38.09050362538873),
  @(45396.99999999999, 18, -0.5589373010535347, 38.56095911678846),
  @(45403.99999999999, 18, -0.4148521179254395, 37.93096876152431),
  @(45410.99999999999, 18, -1.579845201043167, 38.3540447362445),
  @(45417.99999999999, 18, -1.345007966441506, 37.36031527659766),
  @(45424.99999999999, 18, -1.943615987224601, 37.97465373438668),
  @(45438.99999999999, 18, -4.392669045156745, 36.41615305599542),
  @(45445.99999999999, 18, -2.943934781059849, 37.36663733985317),
  @(45452.99999999999, 18, -1.400078664331516, 37.08119882646291),
  @(45487.99999999999, 17, -2.721105401542294, 36.9912677159565),
  @(45494.99999999999, 17, -3.433275280196249, 37.14524257628901),
  @(45501.99999999999, 17, -3.327791591962805, 37.355001503103),
  @(45515.99999999999, 17, -3.608732959998847, 36.60706496544242),
  @(45543.99999999999, 16, -4.86414957370561, 35.7449349970389),
  @(45571.99999999999, 16, -3.453420981411877, 36.44428593319789),
  @(45599.99999999999, 15, -4.488291577872258, 33.89246213633817),
  @(45606.99999999999, 15, -5.868057734316496, 33.86374782711688),
  @(45613.99999999999, 15, -5.162766669770465, 35.44470483553907),
  @(45620.99999999999, 15, -5.059974211077216, 35.54704591711513),
  @(45627.99999999999, 15, -6.331504660329585, 34.85115046376885),
  @(45634.99999999999, 15, -6.795306423675837, 34.4614450795025),
  @(45641.99999999999, 14, -6.100366236110256, 33.55791362056216),
  @(45648.99999999999, 14, -5.298230501610353, 35.03310956804863),
  @(45655.99999999999, 14, -6.742043238187273, 34.47915281802167)
)

$r = 2
foreach ($row in $data) {
  $wsForecast.Cells.Item($r, 1).Value = $row[0]
  $wsForecast.Cells.Item($r, 2).Value = $row[1]
  $wsForecast.Cells.Item($r, 3).Value = $row[2]
  $wsForecast.Cells.Item($r, 4).Value = $row[3]
  $r++
}

# Column A holds the date serials; give it the same custom date/time format
# used on the other two sheets.
$lastRow = $r - 1
$wsForecast.Range("A2:A$lastRow").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Restore the original active sheet (the workbook opened on "Weekly Quantity").
$wsWeekly.Activate()
